$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) { continue }

    $d = $ws.Cells.Item($r, 4).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $f = $ws.Cells.Item($r, 6).Value()

    if ($e -eq 1) {
        $ws.Cells.Item($r, 5).Value = $d
        $ws.Cells.Item($r, 6).Value = $f + $d
    } else {
        $ws.Cells.Item($r, 5).Value = $e - 1
    }
}

Write-Output "Done updating rows"
